$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values ---
$ws.Range("A1").Value = "class name: question"
$ws.Range("B1").Value = "Id: 5"
$ws.Range("C1").Value = "type: generic"
$ws.Range("A2").Value = "associated use cases : 1"
$ws.Range("A3").Value = "operations:"
$ws.Range("B3").Value = "collaborators:"
$ws.Range("A4").Value = "answer"
$ws.Range("B4").Value = "user"
$ws.Range("A5").Value = "create question"
$ws.Range("B5").Value = "user/sprint"
$ws.Range("A6").Value = "report answers"
$ws.Range("B6").Value = "sprint"
$ws.Range("A8").Value = "////////////////////////////////////////////////////////////////////////////////////////////////////////////////////////////////"
$ws.Range("A9").Value = "back"
$ws.Range("A10").Value = "attributes"
$ws.Range("A11").Value = "title"
$ws.Range("B11").Value = "text"
$ws.Range("A12").Value = "description"
$ws.Range("B12").Value = "text"
$ws.Range("A13").Value = "answers"
$ws.Range("B13").Value = "list<(user,text)>"
$ws.Range("A15").Value = "relationships"
$ws.Range("B15").Value = "with"
$ws.Range("A16").Value = "user"
$ws.Range("B16").Value = "user"
$ws.Range("A17").Value = "user"
$ws.Range("B17").Value = "sprint"

# --- Base formatting: Arial 18, centered horizontally & vertically ---
# (applies to the whole used range; row 17 only spans columns A:B)
$ws.Range("A1:C16").Font.Name = "Arial"
$ws.Range("A1:C16").Font.Size = 18
$ws.Range("A1:C16").HorizontalAlignment = -4108
$ws.Range("A1:C16").VerticalAlignment = -4108

$ws.Range("A17:B17").Font.Name = "Arial"
$ws.Range("A17:B17").Font.Size = 18
$ws.Range("A17:B17").HorizontalAlignment = -4108
$ws.Range("A17:B17").VerticalAlignment = -4108

# --- Override to left-aligned for label cells ---
$ws.Range("A1").HorizontalAlignment = -4131
$ws.Range("A2").HorizontalAlignment = -4131
$ws.Range("A3").HorizontalAlignment = -4131
$ws.Range("B3").HorizontalAlignment = -4131
$ws.Range("A9").HorizontalAlignment = -4131
$ws.Range("A10").HorizontalAlignment = -4131
$ws.Range("A15").HorizontalAlignment = -4131
$ws.Range("B15").HorizontalAlignment = -4131

# --- Column widths (characters) ---
$ws.Columns.Item(1).ColumnWidth = 52.59
$ws.Columns.Item(2).ColumnWidth = 29.75
$ws.Columns.Item(3).ColumnWidth = 21.5

# --- Row heights ---
$ws.Range("1:17").RowHeight = 23.25

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection ---
[void]$ws.Range("A6:B6").Select()
